# Updated status as on 10th April
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 57: date corrected from 25-Mar to 27-Mar (40993 -> 40995) ---
$ws.Range("B57").Value = 40995

# --- Row 59: fill in the missing date (was blank) ---
$ws.Range("B59").Value = 40997
$ws.Range("B59").NumberFormat = "d-mmm"

# --- Row 60 ---
$ws.Range("B60").Value = 40998
$ws.Range("B60").NumberFormat = "d-mmm"
$ws.Range("I60").Value = "EKA_METALS_PATCH_0617"
$ws.Range("J60").Value = "Applied to Ref"
$ws.Range("K60").Value = "Released TO QA for UAT"

# --- Row 61 ---
$ws.Range("B61").Value = 41003
$ws.Range("B61").NumberFormat = "d-mmm"
$ws.Range("I61").Value = "EKA_METALS_PATCH_0650"
$ws.Range("J61").Value = "Applied to Ref(METAL_APP_REF)"
$ws.Range("K61").Value = "Released TO QA for UAT"

# --- Row 62 ---
$ws.Range("B62").Value = 41004
$ws.Range("B62").NumberFormat = "d-mmm"
$ws.Range("I62").Value = "EKA_METALS_PATCH_0655"
$ws.Range("J62").Value = "Applied to Ref(METAL_APP_REF)"
$ws.Range("K62").Value = "Released TO QA for UAT"

# --- Highlight the "Applied To Blank" marker column for rows 60-62 ---
$ws.Range("M60").Value = "Applied To Blank"
$ws.Range("M60").Interior.Color = 5296274
$ws.Range("M61").Value = "Applied To Blank"
$ws.Range("M61").Interior.Color = 5296274
$ws.Range("M62").Value = "Applied To Blank"
$ws.Range("M62").Interior.Color = 5296274

# --- Row 63 (newest entry, as of 10th April) ---
$ws.Range("B63").Value = 41009
$ws.Range("B63").NumberFormat = "d-mmm"
$ws.Range("I63").Value = "EKA_METALS_PATCH_0665"
$ws.Range("J63").Value = "Applied to Ref(METAL_APP_REF)"
$ws.Range("M63").Value = "Applied To Blank"
$ws.Range("M63").Interior.Color = 5296274

# --- View state: scroll/zoom/selection to match the latest edits ---
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("I48").Select()
